$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1995
$ws.Range("I32").Value = 1992.6666
$ws.Range("K32").Value = 1992.6666
$ws.Range("M32").Value = -1666.6666
$ws.Range("H40").Value = 1971.5
$ws.Range("I40").Value = 1971.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1971.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1796.5
$ws.Range("H62").Value = 2424.25
$ws.Range("I62").Value = 2424.25
$ws.Range("K62").Value = 2424.25
$ws.Range("M62").Value = -1800.25
$ws.Range("H65").Value = 2424.25
$ws.Range("I65").Value = 2424.25
$ws.Range("K65").Value = 12121.25
$ws.Range("M65").Value = -9001.25
$ws.Range("H74").Value = 6482.731
$ws.Range("I74").Value = 5659.1665
$ws.Range("J74").Value = 7188.643
$ws.Range("K74").Value = 5659.1665
$ws.Range("L74").Value = 7188.643
$ws.Range("M74").Value = -4723.1665
$ws.Range("N74").Value = -9060.643
$ws.Range("H77").Value = 6482.731
$ws.Range("I77").Value = 5659.1665
$ws.Range("J77").Value = 7188.643
$ws.Range("K77").Value = 28295.8325
$ws.Range("L77").Value = 35943.215
$ws.Range("M77").Value = -23615.8325
$ws.Range("N77").Value = -45303.215
$ws.Range("H116").Value = 8271.429
$ws.Range("I116").Value = 9709.182000000001
$ws.Range("K116").Value = 9709.182000000001
$ws.Range("M116").Value = -6267.182000000001
$ws.Range("H129").Value = 2171.5
$ws.Range("I129").Value = 1099
$ws.Range("J129").Value = 4316.5
$ws.Range("K129").Value = 3297
$ws.Range("L129").Value = 12949.5
$ws.Range("M129").Value = 1703
$ws.Range("N129").Value = -22949.5
$ws.Range("H132").Value = 2340.628
$ws.Range("I132").Value = 2162.122
$ws.Range("K132").Value = 6486.366
$ws.Range("M132").Value = -3956.366
$ws.Range("H137").Value = 2958.5833
$ws.Range("I137").Value = 1085.4286
$ws.Range("K137").Value = 3256.2858
$ws.Range("M137").Value = -706.2857999999997
$ws.Range("H138").Value = 3116.4824
$ws.Range("I138").Value = 3421.6191
$ws.Range("J138").Value = 3016.3594
$ws.Range("K138").Value = 10264.8573
$ws.Range("L138").Value = 9049.0782
$ws.Range("M138").Value = -5124.8573
$ws.Range("N138").Value = -19329.0782
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3126.25
$ws.Range("I32").Value = 4028.4285
$ws.Range("K32").Value = 4028.4285
$ws.Range("M32").Value = -3741.4285
$ws.Range("H132").Value = 1917.5536
$ws.Range("I132").Value = 1738.5741
$ws.Range("K132").Value = 5215.7223
$ws.Range("M132").Value = -2685.7223

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 413.23077
$ws.Range("I80").Value = 931.6667
$ws.Range("K80").Value = 931.6667
$ws.Range("M80").Value = 66.33330000000001
$ws.Range("H83").Value = 413.23077
$ws.Range("I83").Value = 931.6667
$ws.Range("K83").Value = 4658.3335
$ws.Range("M83").Value = 333.6665000000003
$ws.Range("H86").Value = 21313.5
$ws.Range("I86").Value = 4300.3335
$ws.Range("K86").Value = 4300.3335
$ws.Range("M86").Value = -3177.3335
$ws.Range("H89").Value = 21313.5
$ws.Range("I89").Value = 4300.3335
$ws.Range("K89").Value = 21501.6675
$ws.Range("M89").Value = -15885.6675
$ws.Range("H94").Value = 5983.769
$ws.Range("I94").Value = 1560.2222
$ws.Range("J94").Value = 15936.75
$ws.Range("K94").Value = 1560.2222
$ws.Range("L94").Value = 15936.75
$ws.Range("M94").Value = -1109.2222
$ws.Range("N94").Value = -16838.75
$ws.Range("H107").Value = 4529.25
$ws.Range("I107").Value = 4529.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4529.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2609.25
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1305.4186
$ws.Range("I31").Value = 623
$ws.Range("J31").Value = 2349.1177
$ws.Range("K31").Value = 623
$ws.Range("L31").Value = 2349.1177
$ws.Range("M31").Value = -328
$ws.Range("N31").Value = -2939.1177
$ws.Range("H34").Value = 1305.4186
$ws.Range("I34").Value = 623
$ws.Range("J34").Value = 2349.1177
$ws.Range("K34").Value = 623
$ws.Range("L34").Value = 2349.1177
$ws.Range("M34").Value = -421
$ws.Range("N34").Value = -2753.1177
$ws.Range("H107").Value = 8476.714
$ws.Range("I107").Value = 17536.334
$ws.Range("J107").Value = 1682
$ws.Range("K107").Value = 17536.334
$ws.Range("L107").Value = 1682
$ws.Range("M107").Value = -15616.334
$ws.Range("N107").Value = -5522
$ws.Range("H132").Value = 1391.9166
$ws.Range("I132").Value = 1431.0526
$ws.Range("J132").Value = 1243.2
$ws.Range("K132").Value = 4293.1578
$ws.Range("L132").Value = 3729.6
$ws.Range("M132").Value = -1763.1578
$ws.Range("N132").Value = -8789.6
$ws.Range("H134").Value = 2623.913
$ws.Range("I134").Value = 2623.913
$ws.Range("K134").Value = 7871.739
$ws.Range("M134").Value = -5336.739
$ws.Range("H141").Value = 46666.332
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 59999.5
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 59999.5
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -70359.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1348.6666
$ws.Range("I60").Value = 1599
$ws.Range("J60").Value = 97
$ws.Range("K60").Value = 4797
$ws.Range("L60").Value = 291
$ws.Range("M60").Value = -4546
$ws.Range("N60").Value = -793
$ws.Range("H68").Value = 1534.5264
$ws.Range("J68").Value = 1534.5264
$ws.Range("L68").Value = 4603.5792
$ws.Range("N68").Value = -6225.5792
$ws.Range("H71").Value = 1534.5264
$ws.Range("J71").Value = 1534.5264
$ws.Range("L71").Value = 13810.7376
$ws.Range("N71").Value = -21922.7376
$ws.Range("H139").Value = 17482.084
$ws.Range("I139").Value = 2039.6666
$ws.Range("J139").Value = 22629.555
$ws.Range("K139").Value = 6118.9998
$ws.Range("L139").Value = 67888.66500000001
$ws.Range("M139").Value = -978.9997999999996
$ws.Range("N139").Value = -78168.66500000001
$ws.Range("H140").Value = 13897927
$ws.Range("I140").Value = 20834274
$ws.Range("K140").Value = 62502822
$ws.Range("M140").Value = -62497642

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2111.0293
$ws.Range("I132").Value = 2123.6897
$ws.Range("J132").Value = 2037.6
$ws.Range("K132").Value = 6371.0691
$ws.Range("L132").Value = 6112.799999999999
$ws.Range("M132").Value = -3841.0691
$ws.Range("N132").Value = -11172.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2733.9
$ws.Range("I7").Value = 2620.5715
$ws.Range("K7").Value = 2620.5715
$ws.Range("M7").Value = -2508.5715
$ws.Range("H16").Value = 2521.9473
$ws.Range("J16").Value = 2847.5
$ws.Range("L16").Value = 2847.5
$ws.Range("N16").Value = -3187.5
$ws.Range("H22").Value = 1399.6364
$ws.Range("J22").Value = 2019.6
$ws.Range("L22").Value = 2019.6
$ws.Range("N22").Value = -2609.6
$ws.Range("H27").Value = 1399.6364
$ws.Range("J27").Value = 2019.6
$ws.Range("L27").Value = 2019.6
$ws.Range("N27").Value = -2233.6
$ws.Range("H68").Value = 2007.75
$ws.Range("I68").Value = 2007.75
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2007.75
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1258.75
$ws.Range("H71").Value = 2007.75
$ws.Range("I71").Value = 2007.75
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 10038.75
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6294.75
$ws.Range("H126").Value = 2733.9
$ws.Range("I126").Value = 2620.5715
$ws.Range("K126").Value = 7861.7145
$ws.Range("M126").Value = -5391.7145
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3963.6667
$ws.Range("I122").Value = 4622.9165
$ws.Range("J122").Value = 1326.6666
$ws.Range("K122").Value = 13868.7495
$ws.Range("L122").Value = 3979.9998
$ws.Range("M122").Value = -11418.7495
$ws.Range("N122").Value = -8879.9998
$ws.Range("H126").Value = 4791.2856
$ws.Range("I126").Value = 2547.25
$ws.Range("J126").Value = 7783.3335
$ws.Range("K126").Value = 7641.75
$ws.Range("L126").Value = 23350.0005
$ws.Range("M126").Value = -5171.75
$ws.Range("N126").Value = -28290.0005
